# This edit reorders the weekly price-report rows (2-20) on the sheet:
# each destination row ends up containing the data that used to live in
# a different (source) row, per the mapping below. Columns A-T are all
# copied together so every field (date, quality, volume, prices, unit,
# origin, $/Kg, Kg/unit) moves as one consistent record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 20
$lastCol  = 20  # column T

# Map of destination row -> source row (both refer to the ORIGINAL,
# pre-edit layout of the sheet).
$perm = @{
    2  = 5
    3  = 6
    4  = 11
    5  = 7
    6  = 12
    7  = 20
    8  = 10
    9  = 15
    10 = 14
    11 = 17
    12 = 3
    13 = 16
    14 = 8
    15 = 18
    16 = 13
    17 = 4
    18 = 2
    19 = 9
    20 = 19
}

# 1) Snapshot every original row's values before we overwrite anything.
$orig = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $row += ,($ws.Cells.Item($r, $c).Value())
    }
    $orig[$r] = $row
}

# 2) Write each destination row using the snapshot of its source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $perm[$r]
    $values = $orig[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
